# Time Planning.xlsx - test edits around the "Bureaucrat" area of the card list.
#
# Net effect (per the target OOXML diff):
#   - C14 ("Council Room") is cut and its text re-appears at D19.
#   - C26 ("Witch") is cut and its text re-appears at D18.
#   - C20 ("Moat") moves one column over to D20.
#   - C28:C33 (Estate, Duchy, Province, Copper, Silver, Gold) each move one
#     column over, from column C to column D, same rows.
#   - The active selection ends up on D20 (it started on C4).
#   - The window/view chrome (workbook xWindow position, and the sheet's
#     scrolled-to topLeftCell) also shifted in the recorded session; we set
#     what this host's object model exposes for that below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell moves (Cut to destination preserves the shared-string text) ---
$ws.Range("C14").Cut($ws.Range("D19"))
$ws.Range("C26").Cut($ws.Range("D18"))
$ws.Range("C20").Cut($ws.Range("D20"))
$ws.Range("C28").Cut($ws.Range("D28"))
$ws.Range("C29").Cut($ws.Range("D29"))
$ws.Range("C30").Cut($ws.Range("D30"))
$ws.Range("C31").Cut($ws.Range("D31"))
$ws.Range("C32").Cut($ws.Range("D32"))
$ws.Range("C33").Cut($ws.Range("D33"))

# --- Selection moves from C4 to D20 ---
$ws.Range("D20").Select()

# --- Window chrome / scroll position (best effort) ---
# Sheet was scrolled so row 3 is the top visible row, and the app window
# shifted horizontally on screen.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Left = 5988
